$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price column (D14:D17)
$ws.Range("D14").Value = 94.12
$ws.Range("D15").Value = 134.55
$ws.Range("D16").Value = 192.4
$ws.Range("D17").Value = 331.5
